# Append new trading log rows (178-191) to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("178", "2026-01-15T01:48:23.859094", "TRADING_ATTEMPT", "SOL", "UNKNOWN", 145.1279744138415, "ATTEMPT", "Attempting trade 1/7"),
    @("179", "2026-01-15T01:48:25.396659", "POSITION_FAILED", "SOL", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 1"),
    @("180", "2026-01-15T01:48:25.456699", "TRADING_ATTEMPT", "ETH", "UNKNOWN", 3327.551283719895, "ATTEMPT", "Attempting trade 2/7"),
    @("181", "2026-01-15T01:48:27.168094", "POSITION_FAILED", "ETH", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 2"),
    @("182", "2026-01-15T01:48:27.226546", "TRADING_ATTEMPT", "ARB", "UNKNOWN", 0.216980312763392, "ATTEMPT", "Attempting trade 3/7"),
    @("183", "2026-01-15T01:48:28.808707", "POSITION_FAILED", "ARB", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 3"),
    @("184", "2026-01-15T01:48:28.867842", "TRADING_ATTEMPT", "AAVE", "UNKNOWN", 176.5633354836947, "ATTEMPT", "Attempting trade 4/7"),
    @("185", "2026-01-15T01:48:30.614313", "POSITION_FAILED", "AAVE", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 4"),
    @("186", "2026-01-15T01:48:30.678715", "TRADING_ATTEMPT", "ADA", "UNKNOWN", 0.4092984782323901, "ATTEMPT", "Attempting trade 5/7"),
    @("187", "2026-01-15T01:48:32.292888", "POSITION_FAILED", "ADA", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 5"),
    @("188", "2026-01-15T01:48:32.353612", "TRADING_ATTEMPT", "ENA", "UNKNOWN", 0.235081410141356, "ATTEMPT", "Attempting trade 6/7"),
    @("189", "2026-01-15T01:48:33.839931", "POSITION_FAILED", "ENA", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 6"),
    @("190", "2026-01-15T01:48:33.900460", "TRADING_ATTEMPT", "DOGE", "UNKNOWN", 0.144971116122147, "ATTEMPT", "Attempting trade 7/7"),
    @("191", "2026-01-15T01:48:35.493078", "POSITION_FAILED", "DOGE", "UNKNOWN", $null, "FAILED", "Trade execution failed for trade 7"),
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    if ($r[5] -ne $null) {
        $ws.Range("E$rowNum").Value = $r[5]
    } else {
        $ws.Range("E$rowNum").Value = ""
    }
    $ws.Range("F$rowNum").Value = ""
    $ws.Range("G$rowNum").Value = ""
    $ws.Range("H$rowNum").Value = ""
    $ws.Range("I$rowNum").Value = ""
    $ws.Range("J$rowNum").Value = ""
    $ws.Range("K$rowNum").Value = $r[6]
    $ws.Range("L$rowNum").Value = $r[7]
}
